$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")
$src = $ws.Range("A100:K100")
$dst = $ws.Range("A100:K102")
$src.AutoFill($dst, 0)
Write-Host "after autofill C101:" $ws.Range("C101").Formula()
Write-Host "after autofill H101:" $ws.Range("H101").Formula()
Write-Host "A101:" $ws.Range("A101").Value2()
